$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 55 entirely; this shifts all rows below it up by one (56->55, ... 96->95)
$ws.Rows("55:55").Delete()

# Update the view: scroll so row 37 is at top-left, and select J47
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Range("J47").Select()
